$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = 36
$ws.Range("F2").Formula = "=(A2+5.15)/0.385"
$ws.Activate()
$ws.Range("F2:F3").Select()
$excel.ActiveWindow.RangeSelection
$ws.Application.ActiveCell

$ws2 = $wb.Worksheets.Item("Arkusz4")
$ws2.Activate()
$ws2.Application.ActiveWindow.ScrollRow = 11
$ws.Activate()
Write-Host "Done"
